$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 106
$ws.Cells.Item(106, 1).Value = "T2L_MD97_2120_d18o_bulloides_SST_from_d18o_bulloides"
$ws.Cells.Item(106, 2).Value = "MD97_2120.Pahnke.2006"
$ws.Cells.Item(106, 3).Value = "Exclude"
$ws.Cells.Item(106, 4).Value = "x"
$ws.Cells.Item(106, 5).Value = "This is a one of those unreviewed 18O records. It has large 2-3C offset with Alkenone and Mg/Ca records from same core, which I would trust more as they have been carefully reviewed by the original authors."

# Row 107
$ws.Cells.Item(107, 1).Value = "ReEnzeIOExA"
$ws.Cells.Item(107, 2).Value = "LakePupuke.Pollen.NewZealand"
$ws.Cells.Item(107, 3).Value = "Exclude"
$ws.Cells.Item(107, 4).Value = "x"
$ws.Cells.Item(107, 5).Value = " I just checked and Darrell appears to have chosen the PLS (~13.5C Holocene average) reconstruction to include in climate12k instead of the MAAT (~14.5C Holocene average), which is favoured by the authors. Maybe he was confusing the the PLS pollen reconstruction with the WAPLS chironomid reconstruction in the same paper. Anyway, the authors definitely favour the MAAT pollen reconstruction so we should use that one."

# Row 108
$ws.Cells.Item(108, 1).Value = "RPZj5YKrFr0"
$ws.Cells.Item(108, 2).Value = "LakePupuke.Pollen.NewZealand"
$ws.Cells.Item(108, 3).Value = "Include and use authors modern temperature of 15C"
$ws.Cells.Item(108, 4).Value = "x"
$ws.Cells.Item(108, 5).Value = "see ReEnzeIOExA.  For information, the paper gives the modern value at the site as 15C, which we could use instead of worldclim"

$ws.Rows.Item(1048576).Delete()

$ws.Range("D108").Select()
